$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(0.3167353711203305,0.2869686923816914,0.2688167699701864,0.2614512282281964,0.2602300919915876,0.2687173079671652,0.306445936546055,0.3814229654958297,0.4371191521956348,0.4625911594059744,0.4722562530874086,0.4701738409460461,0.4633859242391907,0.4592306527138703,0.4354572159095369,0.4209076286605864,0.4125518530815384,0.4097249258378213,0.4224551369059384,0.4653791741896214,0.4935455370996635,0.4785023067315137,0.4217554804644408,0.3610329390560025)
$arr = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt $data.Length; $i++) {
    $arr[$i,0] = $data[$i]
}
$ws.Range("B2:B25").Value = $arr

$data = @(0.03293895998844221,0.02958541580271401,0.02751088255505607,0.02666160096929104,0.02652034306984774,0.02749944461811538,0.03178585670561063,0.04006977794878708,0.04608349490267472,0.04880396840401602,0.04983197765133696,0.04961067403145591,0.04888858687381514,0.04844600470183025,0.04590540062525861,0.04434293716276727,0.04344281848365483,0.04313780768598008,0.04450941209999826,0.04910074031968747,0.05208876626753067,0.05049515655184678,0.04443415450617749,0.03784161320318447)
$arr = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt $data.Length; $i++) {
    $arr[$i,0] = $data[$i]
}
$ws.Range("C2:C25").Value = $arr

$data = @(0.1548485435524185,0.1434574700996265,0.1365594767331757,0.1337724036674999,0.1333110489296345,0.1365217928112088,0.1509007757878891,0.1798746197312724,0.2016562735706842,0.2116773934378031,0.2154886474821183,0.2146670884176842,0.2119906157758535,0.2103533529975437,0.2010036591536775,0.1952969543528766,0.1920251836303919,0.1909192257058834,0.1959033462565927,0.2127763112567465,0.2238999255682756,0.2179541483003931,0.1956291681899316,0.17195110484559)
$arr = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt $data.Length; $i++) {
    $arr[$i,0] = $data[$i]
}
$ws.Range("E2:E25").Value = $arr

$data = @(0.4443680307746263,0.387822817061874,0.3531389305168915,0.3390132514313251,0.336668177824194,0.3529483938344953,0.4248636149813336,0.5661985755041457,0.6702781546542269,0.7176906081379002,0.7356546913071611,0.7317853510981394,0.7191683204515869,0.7114413442032514,0.6671810134426437,0.6400460337125793,0.6244449056556647,0.619163680173358,0.642933953830422,0.7228739723491628,0.7751780083420101,0.7472568307830727,0.6416283278902171,0.5279251897347166)
$arr = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt $data.Length; $i++) {
    $arr[$i,0] = $data[$i]
}
$ws.Range("F2:F25").Value = $arr

$data = @(1.080532715444676,1.06987239300706,1.063942892588656,1.061680795563191,1.061314471107565,1.063911761610228,1.076728780671274,1.106783015149304,1.131911703071353,1.14401582039585,1.148696892808061,1.14768439177891,1.144398976658636,1.14239928503477,1.131134274970918,1.124396438926453,1.120584311144469,1.119304439126608,1.125107136705537,1.145361330567397,1.159167297980048,1.151746515621568,1.124785638583205,1.098120813627418)
$arr = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt $data.Length; $i++) {
    $arr[$i,0] = $data[$i]
}
$ws.Range("G2:G25").Value = $arr

$data = @(1.057445531175318,1.057168693926087,1.057448813117787,1.05767593282954,1.057720462733045,1.057451419014726,1.057256529979639,1.060456253038112,1.065006803016217,1.067558182488881,1.068593799491737,1.068367667328658,1.06764198968969,1.067206545092745,1.064849771190751,1.063527421901114,1.062812129685369,1.062577714622535,1.063663499254005,1.067853251523474,1.070996509829598,1.069281751210042,1.06360183869262,1.059205515511948)
$arr = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt $data.Length; $i++) {
    $arr[$i,0] = $data[$i]
}
$ws.Range("H2:H25").Value = $arr

$data = @(0.272688238731206,0.2421948206752802,0.2235388789508335,0.2159533762494021,0.2146948347070747,0.2234365095121831,0.2621601994152343,0.3386314580823466,0.395148822629011,0.4209347687510956,0.4307102067194535,0.42860440910178,0.4217387825749483,0.4175348007617856,0.3934651730410508,0.3787186375801639,0.3702439688237007,0.3673758211438098,0.3802876920470908,0.4237550892049171,0.452226822546038,0.4370251501899531,0.3795783122087641,0.31788592971958)
$arr = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt $data.Length; $i++) {
    $arr[$i,0] = $data[$i]
}
$ws.Range("K2:K25").Value = $arr

$data = @(0.2282456018431631,0.2073230805264927,0.1945872077282473,0.1894248345044502,0.1885692867944115,0.1945174747041065,0.2210084116627726,0.2738478596039045,0.3132340816596511,0.3312796242979346,0.3381318092960015,0.3366552317309228,0.3318429809027137,0.328897784767392,0.3120573705945802,0.3017594199089544,0.29584839625808,0.2938490934801905,0.302854401811814,0.3332559453444688,0.3532343977230425,0.3425614405411324,0.3023593312212114,0.2594557576664087)
$arr = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt $data.Length; $i++) {
    $arr[$i,0] = $data[$i]
}
$ws.Range("M2:M25").Value = $arr
